# "changed story file internal structure"
# Rebuild the story script on Sheet1 column A with the new scene layout:
# a new intro "Story 0" block is inserted at the top, and the start/end
# markers for each story beat are renumbered (e.g. #SOS1 -> #SOS1-1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "#SOS0-0",
    "스토리 0-1. 1을 눌러 게임 시작",
    "#EOS0-0",
    "#SOS0-1",
    "1 : 진행, 2 : 게임 엔딩 , 3:현제에 머무른다",
    "#EOS0-1",
    "#SOS1-1",
    "스토리 1-1. 앞으로 진행한다",
    "#EOS1-1",
    "#SOS1-2",
    "스토리 1-2. 게임 끝. 다시 처음으로 돌아간다",
    "#EOS1-2",
    "#SOS1-3",
    "스토리 1-3. 현제에 머무른다",
    "#EOS1-3",
    "#SOS2-1",
    "스토리 2-1. 앞으로 진행한다",
    "#EOS2-1",
    "#SOS2-2",
    "스토리 2-2. 게임 끝. 다시 처음으로 돌아간다",
    "#EOS2-2",
    "#SOS2-3",
    "스토리 2-3. 현제에 머무른다",
    "#EOS2-3",
    "#SOS2-4",
    "스토리 2-4. 이전으로 돌아간다",
    "#EOS2-4",
    "#EOF"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value2 = $values[$i]
}

# update the saved selection to match the new sheet state
$ws.Range("A5").Select()
